# Applies the changes described by the commit:
# "add parameter comparison code, change base and v4 parameterizations"
#
# Concretely, on Sheet1:
#  - Delete the entire "Legislature" row (row 12), shifting all rows
#    below it up by one and dropping "Legislature" from the shared
#    strings table.
#  - Adjust a handful of row heights and cell values that were edited
#    by hand alongside the row deletion.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove the "Legislature" row entirely (was row 12).
#    Everything below shifts up by one row automatically.
$ws.Rows.Item(12).Delete()

# 2. Row height tweaks on the rows that remain (post-shift numbering).
$ws.Rows.Item(11).RowHeight = 15    # Bureau of Reclamation
$ws.Rows.Item(12).RowHeight = 15    # Friant Water Authority
$ws.Rows.Item(13).RowHeight = 15    # Friant-Kern Canal
$ws.Rows.Item(14).RowHeight = 29.5  # Irrigation/water districts

# 3. Value updates on "Friant-Kern Canal" row (row 13).
$ws.Range("C13").Value = 1
$ws.Range("D13").Value = 1

# 4. Value updates on "Irrigation/water districts" row (row 14).
$ws.Range("C14").Value = 1
$ws.Range("D14").Value = 1

# 5. Value updates on "Water Rights Division (SWRCB)" row (row 16).
$ws.Range("C16").Value = -1
$ws.Range("D16").Value = -1
$ws.Range("E16").ClearContents()
$ws.Range("F16").ClearContents()

# 6. Update the view's active selection to match the saved workbook.
$ws.Activate()
$ws.Range("C14").Select()
